# CVTest.docx edit: swap the placeholder name/email for the real ones.
#
#   "Victor And"            -> "John Doe"
#   "vicand0502@test.com"   -> "Johndoe0502@test.com", but split so the
#                              <w:br/> before it, "Johndoe" and
#                              "0502@test.com" each live in their own run
#                              (matching the target markup).

$d = $word.ActiveDocument

# 1) Name: "Victor And" -> "John Doe"
$d.Content.Find.Execute("Victor And", $true, $false, $false, $false, $false, $true, 1, $false, "John Doe", 2)

# 2) Email: swap the text in place first (keeps it in the same run as the
#    preceding <w:br/> for now; we split that apart next).
$d.Content.Find.Execute("vicand0502@test.com", $true, $false, $false, $false, $false, $true, 1, $false, "Johndoe0502@test.com", 2)

# 3) Re-locate the new text so we can compute exact split points.
$seg = $d.Content
$seg.Find.Execute("Johndoe0502@test.com")
$segStart = $seg.Start
$segEnd   = $seg.End
$brPos    = $segStart - 1   # the <w:br/> immediately precedes this text
$mid      = $segStart + 7   # length of "Johndoe"

# 3a) Isolate the <w:br/> into its own run, separate from "John Doe" before
#     it and "Johndoe" after it. Toggling a character-formatting property on
#     and back off is the standard COM trick to force Word to break the run
#     at these boundaries without changing the visible formatting.
$rBreak = $d.Range($brPos, $segStart)
$rBreak.Bold = $true
$rBreak.Bold = $false

# 3b) Split "Johndoe" from "0502@test.com" into two separate runs.
$rTail = $d.Range($mid, $segEnd)
$rTail.Bold = $true
$rTail.Bold = $false

Write-Output $d.Paragraphs(1).Range.Text
